# Functional Programming Update to Main Program
# Applies the "route_3" columns to Routes, tweaks Demands figures, and
# fleshes out the fuel/electric reference table on Possible_Updates.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Routes": add route_3_index / route_3_mean / route_3_std columns
# (M, N, O) mirroring route_1's shape, and correct refill_std for row 3.
# ---------------------------------------------------------------------
$wsRoutes = $wb.Worksheets.Item("Routes")

$wsRoutes.Range("F3").Value = 0.25

$wsRoutes.Range("M1").Value = "route_3_index"
$wsRoutes.Range("N1").Value = "route_3_mean"
$wsRoutes.Range("O1").Value = "route_3_std"

$wsRoutes.Range("M2").Value = 1
$wsRoutes.Range("N2").Value = 5
$wsRoutes.Range("O2").Value = 1

$wsRoutes.Range("M3").Value = 0
$wsRoutes.Range("N3").Value = 1
$wsRoutes.Range("O3").Value = 0.5

$wsRoutes.Range("M4").Value = 1
$wsRoutes.Range("N4").Value = 7
$wsRoutes.Range("O4").Value = 2

$wsRoutes.Range("M5").Value = 0
$wsRoutes.Range("N5").Value = 0.75
$wsRoutes.Range("O5").Value = 1.5

$wsRoutes.Range("M6").Value = 1
$wsRoutes.Range("N6").Value = 5
$wsRoutes.Range("O6").Value = 1.5

# New columns mirror the width of the route_1 / route_2 columns, which
# already have an identical best-fit shape (same header/value lengths).
$wsRoutes.Columns.Item(13).ColumnWidth = $wsRoutes.Columns.Item(7).ColumnWidth
$wsRoutes.Columns.Item(14).ColumnWidth = $wsRoutes.Columns.Item(8).ColumnWidth
$wsRoutes.Columns.Item(15).ColumnWidth = $wsRoutes.Columns.Item(9).ColumnWidth

# ---------------------------------------------------------------------
# Sheet "Demands": update a few demand figures
# ---------------------------------------------------------------------
$wsDemands = $wb.Worksheets.Item("Demands")

$wsDemands.Range("B2").Value = 3
$wsDemands.Range("C2").Value = 2
$wsDemands.Range("F3").Value = 25
$wsDemands.Range("B4").Value = 3
$wsDemands.Range("F4").Value = 25

# ---------------------------------------------------------------------
# Sheet "Possible_Updates": fuel / electric reference table
# ---------------------------------------------------------------------
$wsUpdates = $wb.Worksheets.Item("Possible_Updates")

$wsUpdates.Range("B6").Value = "Gallons per hour"
$wsUpdates.Range("C6").Value = "lts per minute"

$wsUpdates.Range("A7").Value = "Fuel Running"
$wsUpdates.Range("E7").Formula = "=8.9/100"
$wsUpdates.Range("F7").Value = 48.28
$wsUpdates.Range("F8").Formula = "=E7*F7"
$wsUpdates.Range("C7").Formula = "=F8/60"
$wsUpdates.Range("G7").Value = "https://www.energy.gov/eere/vehicles/fact-861-february-23-2015-idle-fuel-consumption-selected-gasoline-and-diesel-vehicles"

$wsUpdates.Range("A8").Value = "Fuel Idle"
$wsUpdates.Range("B8").Value = 0.97
$wsUpdates.Range("C8").Formula = "=3.78*B8/60"
$wsUpdates.Range("G8").Value = "https://www.google.com/search?q=average+fuel+consumtion+of+transit+bus+per+hour&rlz=1C1GCEA_enUS970US970&oq=average+fuel+consumtion+of+transit+bus+per+hour&aqs=chrome..69i57j33i10i22i29i30.16846j0j9&sourceid=chrome&ie=UTF-8"

$wsUpdates.Range("A9").Value = "ReFueling"
$wsUpdates.Range("C9").Value = "30-34"
$wsUpdates.Range("G9").Value = "https://www.quora.com/What-is-the-flow-rate-of-gasoline-station-fuel-dispensers-in-liters-per-second"

$wsUpdates.Range("A10").Value = "Electric"

$wsUpdates.Columns.Item(1).ColumnWidth = 10.944010416666666
$wsUpdates.Columns.Item(2).ColumnWidth = 13.608072916666666
$wsUpdates.Columns.Item(3).ColumnWidth = 11.276041666666666

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------
$wsDemands.Range("G2").Select() | Out-Null
$wsUpdates.Range("C7").Select() | Out-Null

$wsRoutes.Activate() | Out-Null
$wsRoutes.Range("O3").Select() | Out-Null
